# Update cryptos list (prices + 1h volume %) per Mon Oct 23 11:54:32 UTC 2023 refresh
# Rows 13/14 (Chainlink <-> WrappedEther) also swap order/data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "30.656.99"
$ws.Cells.Item(2, 5).Value = "  +2.52%  "

$ws.Cells.Item(3, 4).Value = "1.676.44"
$ws.Cells.Item(3, 5).Value = "  +2.92%  "

$ws.Cells.Item(4, 5).Value = "  -0.44%  "

$ws.Cells.Item(5, 4).Value = "'219.83"
$ws.Cells.Item(5, 5).Value = "  +2.69%  "

$ws.Cells.Item(6, 5).Value = "  +2.25%  "

$ws.Cells.Item(7, 5).Value = "  -0.33%  "

$ws.Cells.Item(8, 4).Value = "'29.66"
$ws.Cells.Item(8, 5).Value = "  +4.90%  "

$ws.Cells.Item(9, 5).Value = "  +3.01%  "

$ws.Cells.Item(10, 4).Value = "'0.0650"
$ws.Cells.Item(10, 5).Value = "  +7.14%  "

$ws.Cells.Item(11, 5).Value = "  -0.12%  "

$ws.Cells.Item(12, 4).Value = "1.916.57"
$ws.Cells.Item(12, 5).Value = "  +2.87%  "

$ws.Cells.Item(13, 2).Value = "WrappedEther"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(13, 4).Value = "1.718.94"
$ws.Cells.Item(13, 5).Value = "  +5.50%  "

$ws.Cells.Item(14, 2).Value = "Chainlink"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(14, 4).Value = "'10.27"
$ws.Cells.Item(14, 5).Value = "  +13.31%  "

$ws.Cells.Item(15, 5).Value = "  +9.51%  "

$ws.Cells.Item(16, 5).Value = "  +4.83%  "

$ws.Cells.Item(17, 4).Value = "30.658.89"
$ws.Cells.Item(17, 5).Value = "  +2.43%  "

$ws.Cells.Item(18, 4).Value = "'66.41"
$ws.Cells.Item(18, 5).Value = "  +3.86%  "

$ws.Cells.Item(19, 4).Value = "'243.93"
$ws.Cells.Item(19, 5).Value = "  +1.18%  "

$ws.Cells.Item(20, 4).Value = "0.0₃0726"
$ws.Cells.Item(20, 5).Value = "  +3.72%  "

$ws.Cells.Item(21, 5).Value = "  -0.32%  "

$ws.Cells.Item(22, 4).Value = "'4.27"
$ws.Cells.Item(22, 5).Value = "  +3.95%  "

$ws.Cells.Item(23, 4).Value = "'10.01"
$ws.Cells.Item(23, 5).Value = "  +3.17%  "

$ws.Cells.Item(24, 5).Value = "  +0.49%  "

$ws.Cells.Item(25, 4).Value = "'158.69"
$ws.Cells.Item(25, 5).Value = "  -0.07%  "

$ws.Cells.Item(26, 4).Value = "'15.89"
$ws.Cells.Item(26, 5).Value = "  +2.65%  "

$ws.Cells.Item(27, 5).Value = "  +2.93%  "

$ws.Cells.Item(28, 5).Value = "  +1.82%  "

$ws.Cells.Item(29, 5).Value = "  -0.37%  "

$ws.Cells.Item(30, 5).Value = "  +2.24%  "

$ws.Cells.Item(31, 4).Value = "'1.14"
$ws.Cells.Item(31, 5).Value = "  +3.61%  "

$ws.Cells.Item(32, 4).Value = "'3.46"
$ws.Cells.Item(32, 5).Value = "  +3.22%  "

$ws.Cells.Item(33, 4).Value = "'3.29"
$ws.Cells.Item(33, 5).Value = "  +4.18%  "

$ws.Cells.Item(34, 4).Value = "1.487.81"
$ws.Cells.Item(34, 5).Value = "  +4.57%  "

$ws.Cells.Item(35, 5).Value = "  +7.80%  "

$ws.Cells.Item(36, 4).Value = "'85.17"
$ws.Cells.Item(36, 5).Value = "  +13.44%  "

$ws.Cells.Item(37, 5).Value = "  -0.35%  "

$ws.Cells.Item(38, 4).Value = "'0.601"
$ws.Cells.Item(38, 5).Value = "  +9.38%  "

$ws.Cells.Item(39, 4).Value = "'0.0179"
$ws.Cells.Item(39, 5).Value = "  +5.98%  "

$ws.Cells.Item(40, 5).Value = "  -3.28%  "

$ws.Cells.Item(41, 5).Value = "  -0.40%  "

$ws.Cells.Item(42, 4).Value = "'0.840"
$ws.Cells.Item(42, 5).Value = "  +2.00%  "

$ws.Cells.Item(43, 5).Value = "  +1.74%  "

$ws.Cells.Item(44, 5).Value = "  -0.27%  "

$ws.Cells.Item(45, 5).Value = "  +0.09%  "

$ws.Cells.Item(46, 5).Value = "  -0.30%  "

$ws.Cells.Item(47, 4).Value = "'51.46"
$ws.Cells.Item(47, 5).Value = "  -1.63%  "

$ws.Cells.Item(48, 5).Value = "  +3.37%  "

$ws.Cells.Item(49, 4).Value = "1.808.15"
$ws.Cells.Item(49, 5).Value = "  +2.12%  "

$ws.Cells.Item(50, 4).Value = "'94.77"
$ws.Cells.Item(50, 5).Value = "  +5.03%  "

$ws.Cells.Item(51, 4).Value = "0.0⁦0113"
$ws.Cells.Item(51, 5).Value = "  -0.43%  "
